$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 69 (the existing row 69 and everything below
# shifts down by one). The new row inherits formatting from row 68 above it,
# which already carries the "mini header" look (B/C bold key-style, D/E
# blank) used for standalone localization keys.
$ws.Rows.Item(69).Insert()

# Populate the new row with the new localization key/value pair.
$ws.Range("B69").Value = "relatedQuestion"
$ws.Range("C69").Value = "Related question"
